$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1068.4314
$ws.Cells.Item(15, 9).Value = 1068.4314
$ws.Cells.Item(15, 11).Value = 3205.2942
$ws.Cells.Item(15, 13).Value = -3036.2942

$ws.Cells.Item(93, 8).Value = 79999.5
$ws.Cells.Item(93, 9).Value = 80000
$ws.Cells.Item(93, 10).Value = 79999
$ws.Cells.Item(93, 11).Value = 80000
$ws.Cells.Item(93, 12).Value = 79999
$ws.Cells.Item(93, 13).Value = -77504
$ws.Cells.Item(93, 14).Value = -84991

$ws.Cells.Item(96, 8).Value = 2949.889
$ws.Cells.Item(96, 9).Value = 2269.8
$ws.Cells.Item(96, 11).Value = 6809.400000000001
$ws.Cells.Item(96, 13).Value = -5436.400000000001

$ws.Cells.Item(116, 8).Value = 10000
$ws.Cells.Item(116, 10).Value = 10000
$ws.Cells.Item(116, 12).Value = 10000
$ws.Cells.Item(116, 14).Value = -16884

$ws.Cells.Item(137, 8).Value = 2296.5386
$ws.Cells.Item(137, 10).Value = 3576.2222
$ws.Cells.Item(137, 12).Value = 10728.6666
$ws.Cells.Item(137, 14).Value = -15828.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1927.75
$ws.Cells.Item(2, 9).Value = 903.6667
$ws.Cells.Item(2, 11).Value = 903.6667
$ws.Cells.Item(2, 13).Value = -790.6667

$ws.Cells.Item(32, 8).Value = 5469.1636
$ws.Cells.Item(32, 9).Value = 2949.186
$ws.Cells.Item(32, 11).Value = 2949.186
$ws.Cells.Item(32, 13).Value = -2662.186

$ws.Cells.Item(38, 8).Value = 10000
$ws.Cells.Item(38, 9).Value = 10000
$ws.Cells.Item(38, 11).Value = 10000
$ws.Cells.Item(38, 13).Value = -9533

$ws.Cells.Item(45, 8).Value = 1923.3334
$ws.Cells.Item(45, 9).Value = 1898
$ws.Cells.Item(45, 11).Value = 1898
$ws.Cells.Item(45, 13).Value = -1521

$ws.Cells.Item(92, 8).Value = 86498
$ws.Cells.Item(92, 9).Value = 80000
$ws.Cells.Item(92, 10).Value = 89747
$ws.Cells.Item(92, 11).Value = 80000
$ws.Cells.Item(92, 12).Value = 89747
$ws.Cells.Item(92, 13).Value = -77504
$ws.Cells.Item(92, 14).Value = -94739

$ws.Cells.Item(116, 8).Value = 1927.75
$ws.Cells.Item(116, 9).Value = 903.6667
$ws.Cells.Item(116, 11).Value = 903.6667
$ws.Cells.Item(116, 13).Value = 1390.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1927.75
$ws.Cells.Item(3, 9).Value = 903.6667
$ws.Cells.Item(3, 11).Value = 903.6667
$ws.Cells.Item(3, 13).Value = -789.6667

$ws.Cells.Item(92, 8).Value = 32500
$ws.Cells.Item(92, 9).Value = 40000
$ws.Cells.Item(92, 10).Value = 25000
$ws.Cells.Item(92, 11).Value = 40000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 13).Value = -37504
$ws.Cells.Item(92, 14).Value = -29992

$ws.Cells.Item(107, 8).Value = 1100
$ws.Cells.Item(107, 9).Value = 1088.8889
$ws.Cells.Item(107, 11).Value = 1088.8889
$ws.Cells.Item(107, 13).Value = 831.1111000000001

$ws.Cells.Item(134, 8).Value = 2760.1177
$ws.Cells.Item(134, 9).Value = 2546
$ws.Cells.Item(134, 11).Value = 7638
$ws.Cells.Item(134, 13).Value = -5103

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 502.5
$ws.Cells.Item(22, 10).Value = 406
$ws.Cells.Item(22, 12).Value = 406
$ws.Cells.Item(22, 14).Value = -1106

$ws.Cells.Item(31, 8).Value = 5679.263
$ws.Cells.Item(31, 9).Value = 4026.7334
$ws.Cells.Item(31, 10).Value = 6757
$ws.Cells.Item(31, 11).Value = 4026.7334
$ws.Cells.Item(31, 12).Value = 6757
$ws.Cells.Item(31, 13).Value = -3731.7334
$ws.Cells.Item(31, 14).Value = -7347

$ws.Cells.Item(34, 8).Value = 5679.263
$ws.Cells.Item(34, 9).Value = 4026.7334
$ws.Cells.Item(34, 10).Value = 6757
$ws.Cells.Item(34, 11).Value = 4026.7334
$ws.Cells.Item(34, 12).Value = 6757
$ws.Cells.Item(34, 13).Value = -3824.7334
$ws.Cells.Item(34, 14).Value = -7161

$ws.Cells.Item(35, 8).Value = 13997.5
$ws.Cells.Item(35, 9).Value = 19995
$ws.Cells.Item(35, 10).Value = 8000
$ws.Cells.Item(35, 11).Value = 19995
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = -19701
$ws.Cells.Item(35, 14).Value = -8588

$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 13).ClearContents()

$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 13).ClearContents()

$ws.Cells.Item(58, 8).Value = 3529.4546
$ws.Cells.Item(58, 9).Value = 1289.7142
$ws.Cells.Item(58, 10).Value = 4574.6665
$ws.Cells.Item(58, 11).Value = 1289.7142
$ws.Cells.Item(58, 12).Value = 4574.6665
$ws.Cells.Item(58, 13).Value = -1086.7142
$ws.Cells.Item(58, 14).Value = -4980.6665

$ws.Cells.Item(99, 8).Value = 11175.2
$ws.Cells.Item(99, 9).Value = 6383.357
$ws.Cells.Item(99, 11).Value = 6383.357
$ws.Cells.Item(99, 13).Value = -4885.357

$ws.Cells.Item(122, 8).Value = 2246.7827
$ws.Cells.Item(122, 9).Value = 2991.1538
$ws.Cells.Item(122, 10).Value = 1279.1
$ws.Cells.Item(122, 11).Value = 8973.4614
$ws.Cells.Item(122, 12).Value = 3837.3
$ws.Cells.Item(122, 13).Value = -6523.4614
$ws.Cells.Item(122, 14).Value = -8737.299999999999

$ws.Cells.Item(126, 8).Value = 11175.2
$ws.Cells.Item(126, 9).Value = 6383.357
$ws.Cells.Item(126, 11).Value = 19150.071
$ws.Cells.Item(126, 13).Value = -16680.071

$ws.Cells.Item(132, 8).Value = 1968.1428
$ws.Cells.Item(132, 9).Value = 1731.9412
$ws.Cells.Item(132, 11).Value = 5195.8236
$ws.Cells.Item(132, 13).Value = -2665.8236

$ws.Cells.Item(136, 8).Value = 3529.4546
$ws.Cells.Item(136, 9).Value = 1289.7142
$ws.Cells.Item(136, 10).Value = 4574.6665
$ws.Cells.Item(136, 11).Value = 3869.1426
$ws.Cells.Item(136, 12).Value = 13723.9995
$ws.Cells.Item(136, 13).Value = -1319.1426
$ws.Cells.Item(136, 14).Value = -18823.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 95
$ws.Cells.Item(8, 9).Value = 95
$ws.Cells.Item(8, 11).Value = 285
$ws.Cells.Item(8, 13).Value = -146

$ws.Cells.Item(46, 8).Value = 5000500
$ws.Cells.Item(46, 10).Value = 5000500
$ws.Cells.Item(46, 12).Value = 15001500
$ws.Cells.Item(46, 14).Value = -15001682

$ws.Cells.Item(82, 8).Value = 8716.666999999999
$ws.Cells.Item(82, 10).Value = 8460
$ws.Cells.Item(82, 12).Value = 25380
$ws.Cells.Item(82, 14).Value = -26192

$ws.Cells.Item(85, 8).Value = 8716.666999999999
$ws.Cells.Item(85, 10).Value = 8460
$ws.Cells.Item(85, 12).Value = 25380
$ws.Cells.Item(85, 14).Value = -28188

$ws.Cells.Item(131, 8).Value = 1037.8889
$ws.Cells.Item(131, 9).Value = 198.83333
$ws.Cells.Item(131, 10).Value = 1457.4166
$ws.Cells.Item(131, 11).Value = 596.49999
$ws.Cells.Item(131, 12).Value = 4372.2498
$ws.Cells.Item(131, 13).Value = 4443.50001
$ws.Cells.Item(131, 14).Value = -14452.2498

$ws.Cells.Item(132, 8).Value = 11196.8
$ws.Cells.Item(132, 9).Value = 10996
$ws.Cells.Item(132, 11).Value = 98964
$ws.Cells.Item(132, 13).Value = -96434

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2135.1082
$ws.Cells.Item(102, 9).Value = 943.6429000000001
$ws.Cells.Item(102, 11).Value = 943.6429000000001
$ws.Cells.Item(102, 13).Value = 678.3570999999999

$ws.Cells.Item(122, 8).Value = 74142.36
$ws.Cells.Item(122, 9).Value = 2585.7778
$ws.Cells.Item(122, 11).Value = 7757.3334
$ws.Cells.Item(122, 13).Value = -5307.3334

$ws.Cells.Item(126, 8).Value = 4202.75
$ws.Cells.Item(126, 9).Value = 3152
$ws.Cells.Item(126, 11).Value = 9456
$ws.Cells.Item(126, 13).Value = -6986

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 9285.714
$ws.Cells.Item(2, 9).Value = 9461.538
$ws.Cells.Item(2, 10).Value = 7000
$ws.Cells.Item(2, 11).Value = 9461.538
$ws.Cells.Item(2, 12).Value = 7000
$ws.Cells.Item(2, 13).Value = -9349.538
$ws.Cells.Item(2, 14).Value = -7224

$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()

$ws.Cells.Item(76, 8).Value = 34750
$ws.Cells.Item(76, 10).Value = 34750
$ws.Cells.Item(76, 12).Value = 34750
$ws.Cells.Item(76, 14).Value = -35426

$ws.Cells.Item(79, 8).Value = 34750
$ws.Cells.Item(79, 10).Value = 34750
$ws.Cells.Item(79, 12).Value = 34750
$ws.Cells.Item(79, 14).Value = -37090

$ws.Cells.Item(108, 8).Value = 99333.336
$ws.Cells.Item(108, 10).Value = 19000
$ws.Cells.Item(108, 12).Value = 19000
$ws.Cells.Item(108, 14).Value = -26680

$ws.Cells.Item(122, 8).Value = 6233.3125
$ws.Cells.Item(122, 9).Value = 7411.375
$ws.Cells.Item(122, 11).Value = 22234.125
$ws.Cells.Item(122, 13).Value = -19784.125

$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 3513.3635
$ws.Cells.Item(132, 9).Value = 2733.8635
$ws.Cells.Item(132, 10).Value = 5072.364
$ws.Cells.Item(132, 11).Value = 8201.5905
$ws.Cells.Item(132, 12).Value = 15217.092
$ws.Cells.Item(132, 13).Value = -5671.5905
$ws.Cells.Item(132, 14).Value = -20277.092

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8453.111000000001
$ws.Cells.Item(62, 9).Value = 8299
$ws.Cells.Item(62, 10).Value = 8497.143
$ws.Cells.Item(62, 11).Value = 8299
$ws.Cells.Item(62, 12).Value = 8497.143
$ws.Cells.Item(62, 13).Value = -7675
$ws.Cells.Item(62, 14).Value = -9745.143

$ws.Cells.Item(65, 8).Value = 8453.111000000001
$ws.Cells.Item(65, 9).Value = 8299
$ws.Cells.Item(65, 10).Value = 8497.143
$ws.Cells.Item(65, 11).Value = 41495
$ws.Cells.Item(65, 12).Value = 42485.715
$ws.Cells.Item(65, 13).Value = -38375
$ws.Cells.Item(65, 14).Value = -48725.715
